$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.091.70'
$ws.Range('E2').Value = '  +1.49%  '

# Row 3
$ws.Range('D3').Value = '2.046.37'
$ws.Range('E3').Value = '  +0.80%  '

# Row 4
$ws.Range('E4').Value = '  +0.25%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '231.22'
$ws.Range('E5').Value = '  -0.43%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  +2.55%  '

# Row 7
$ws.Range('E7').Value = '  +0.13%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '57.10'
$ws.Range('E8').Value = '  +3.76%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.379'
$ws.Range('E9').Value = '  +2.46%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '57.65'
$ws.Range('E10').Value = '  +0.91%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0753'
$ws.Range('E11').Value = '  +0.87%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.101'
$ws.Range('E12').Value = '  +0.82%  '

# Row 13
$ws.Range('D13').Value = '2.350.34'
$ws.Range('E13').Value = '  +1.20%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.22'
$ws.Range('E14').Value = '  -0.35%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.62'
$ws.Range('E15').Value = '  +2.31%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.767'
$ws.Range('E16').Value = '  +0.89%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.12'
$ws.Range('E17').Value = '  -0.09%  '

# Row 18
$ws.Range('D18').Value = '2.053.49'
$ws.Range('E18').Value = '  +1.19%  '

# Row 19
$ws.Range('D19').Value = '37.068.24'
$ws.Range('E19').Value = '  +1.19%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.25'
$ws.Range('E20').Value = '  +13.94%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '68.80'
$ws.Range('E21').Value = '  +1.91%  '

# Row 22
$ws.Range('D22').Value = '0.0₃0804'
$ws.Range('E22').Value = '  +1.03%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '223.79'
$ws.Range('E23').Value = '  +1.34%  '

# Row 24
$ws.Range('E24').Value = '  -0.02%  '

# Row 25
$ws.Range('E25').Value = '  +1.41%  '

# Row 26
$ws.Range('E26').Value = '  -0.34%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.12'
$ws.Range('E27').Value = '  +1.57%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.44'
$ws.Range('E28').Value = '  +6.48%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.71'
$ws.Range('E29').Value = '  +0.48%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '18.91'
$ws.Range('E30').Value = '  +0.25%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.124'
$ws.Range('E31').Value = '  -3.67%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.116'
$ws.Range('E32').Value = '  -1.22%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.42'
$ws.Range('E33').Value = '  +0.98%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0608'
$ws.Range('E34').Value = '  +1.02%  '

# Row 35
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.52'
$ws.Range('E35').Value = '  +1.96%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.50'
$ws.Range('E36').Value = '  +5.38%  '

# Row 37
$ws.Range('E37').Value = '  +0.29%  '

# Row 38
$ws.Range('E38').Value = '  -0.85%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.24'
$ws.Range('E39').Value = '  -1.36%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.70'
$ws.Range('E40').Value = '  -2.08%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.57'
$ws.Range('E41').Value = '  +9.39%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.94'
$ws.Range('E42').Value = '  +1.01%  '

# Row 43
$ws.Range('D43').Value = '1.475.84'
$ws.Range('E43').Value = '  +0.44%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '95.51'
$ws.Range('E44').Value = '  +3.06%  '

# Row 45
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0924'
$ws.Range('E45').Value = '  -0.17%  '

# Row 46
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.16'
$ws.Range('E46').Value = '  +3.16%  '

# Row 47
$ws.Range('E47').Value = '  +2.81%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.01'
$ws.Range('E48').Value = '  +1.01%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '15.14'
$ws.Range('E49').Value = '  -2.91%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.09'
$ws.Range('E50').Value = '  +3.33%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.93'
$ws.Range('E51').Value = '  +1.15%  '
